$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as plain text in the source data
# (e.g. "20.30", "9.19", "36.979.03"). Excel auto-converts numeric-looking
# strings assigned via .Value into real numbers, which would silently drop
# trailing zeros / change the stored type from Text to Number. Forcing the
# cell number format to Text ("@") before assignment keeps the literal string.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.979.03"
$ws.Range("E2").Value = "  -0.91%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.011.36"
$ws.Range("E3").Value = "  -1.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.03"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.57"
$ws.Range("E8").Value = "  -4.00%  "

$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  +1.92%  "

$ws.Range("E11").Value = "  -3.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.307.77"
$ws.Range("E12").Value = "  -1.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.25"
$ws.Range("E13").Value = "  -2.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.30"
$ws.Range("E14").Value = "  -0.76%  "

$ws.Range("E15").Value = "  -1.47%  "

$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.011.28"
$ws.Range("E17").Value = "  -1.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.874.18"
$ws.Range("E18").Value = "  -1.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  +3.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.73"
$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.80"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  +2.86%  "

$ws.Range("E25").Value = "  -5.58%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.19"
$ws.Range("E26").Value = "  -2.84%  "

$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.74"
$ws.Range("E27").Value = "  -1.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.126"
$ws.Range("E28").Value = "  -1.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.68"
$ws.Range("E29").Value = "  -1.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  -1.24%  "

$ws.Range("E31").Value = "  -3.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.49"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("E34").Value = "  -3.89%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.35"
$ws.Range("E35").Value = "  -3.75%  "

$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("E38").Value = "  -3.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.33"
$ws.Range("E39").Value = "  +2.38%  "

$ws.Range("E40").Value = "  -3.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.478.30"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.78"
$ws.Range("E42").Value = "  +2.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.15"
$ws.Range("E43").Value = "  -2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0918"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.78"
$ws.Range("E45").Value = "  -3.99%  "

$ws.Range("E46").Value = "  -4.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("E47").Value = "  +0.71%  "

$ws.Range("E48").Value = "  -1.89%  "

$ws.Range("E49").Value = "  +0.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.196.79"
$ws.Range("E50").Value = "  -1.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.61"
$ws.Range("E51").Value = "  -10.74%  "
